# Atualizado por script em 25-10-2023 12:32
#
# Swap the match-data columns (F:V) between the row pairs that were
# re-ordered by the upstream scrape (the A:E "index / league / date"
# columns stay pinned to the row, only the match content moved rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, $rowA, $rowB)
    $rangeA = $ws.Range("F$rowA`:V$rowA")
    $rangeB = $ws.Range("F$rowB`:V$rowB")
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

Swap-RowData $ws 67 68
Swap-RowData $ws 73 74
Swap-RowData $ws 82 83

# Append the four new match rows (108-111) that were scraped since the
# previous run. Copy A/E formatting from the last existing row (107) so
# the new rows match the sheet's styling.

$ws.Range("A107:E107").Copy()
$ws.Range("A108:E111").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row=108; A=107; E=45224.625; F="Jadran Dekani";   G=1; H="NK Bistrica";  I=4;
       J=2.31; K="24/10/2023 02:12"; L=2.05; M="25/10/2023 14:57";
       N=3.11; O="24/10/2023 02:12"; P=3.32; Q="25/10/2023 14:58";
       R=2.68; S="24/10/2023 02:12"; T=3.37; U="25/10/2023 14:58";
       V="https://www.betexplorer.com/football/slovenia/2-snl/jadran-dekani-bistrica/WnwnbggE/" },
    @{ Row=109; A=108; E=45224.625; F="Beltinci";         G=1; H="Triglav";      I=2;
       J=2.13; K="08/08/2023 04:42"; L=1.49; M="25/10/2023 14:59";
       N=3.25; O="08/08/2023 04:42"; P=4.51; Q="25/10/2023 14:59";
       R=2.86; S="08/08/2023 04:42"; T=5.21; U="25/10/2023 14:59";
       V="https://www.betexplorer.com/football/slovenia/2-snl/beltinci-triglav/xUgqvY1r/" },
    @{ Row=110; A=109; E=45224.625; F="Dravinja";         G=3; H="NK Krka";      I=0;
       J=3.2;  K="24/10/2023 02:12"; L=2.88; M="25/10/2023 14:57";
       N=3.36; O="24/10/2023 02:12"; P=3.39; Q="25/10/2023 14:57";
       R=1.93; S="24/10/2023 02:12"; T=2.26; U="25/10/2023 14:57";
       V="https://www.betexplorer.com/football/slovenia/2-snl/dravinja-nk-krka/AwxjcD8K/" },
    @{ Row=111; A=110; E=45224.625; F="Tolmin";           G=3; H="Tabor Sezana"; I=1;
       J=1.69; K="24/10/2023 02:12"; L=1.96; M="25/10/2023 14:50";
       N=3.65; O="24/10/2023 02:12"; P=3.8;  Q="25/10/2023 14:50";
       R=3.78; S="24/10/2023 02:12"; T=3.2;  U="25/10/2023 14:50";
       V="https://www.betexplorer.com/football/slovenia/2-snl/tolmin-tabor-sezana/GKhmwhHl/" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value2 = $r.A
    $ws.Range("B$row").Value2 = "slovenia"
    $ws.Range("C$row").Value2 = "2-snl"
    $ws.Range("D$row").Value2 = "2023-2024"
    $ws.Range("E$row").Value2 = $r.E
    $ws.Range("F$row").Value2 = $r.F
    $ws.Range("G$row").Value2 = $r.G
    $ws.Range("H$row").Value2 = $r.H
    $ws.Range("I$row").Value2 = $r.I
    $ws.Range("J$row").Value2 = $r.J
    $ws.Range("K$row").Value2 = $r.K
    $ws.Range("L$row").Value2 = $r.L
    $ws.Range("M$row").Value2 = $r.M
    $ws.Range("N$row").Value2 = $r.N
    $ws.Range("O$row").Value2 = $r.O
    $ws.Range("P$row").Value2 = $r.P
    $ws.Range("Q$row").Value2 = $r.Q
    $ws.Range("R$row").Value2 = $r.R
    $ws.Range("S$row").Value2 = $r.S
    $ws.Range("T$row").Value2 = $r.T
    $ws.Range("U$row").Value2 = $r.U
    $ws.Range("V$row").Value2 = $r.V
}
